$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.708.98"
$ws.Range("E2").Value = "  +3.13%  "
$ws.Range("D3").Value = "3.281.75"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "3.273.69"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.575"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000270"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "687.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.21%  "
$ws.Range("D15").Value = "3.802.88"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "67.787.44"
$ws.Range("E17").Value = "  +3.22%  "
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "3.278.94"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.893"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "582.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.40%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "3.860.70"
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.104"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "32.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").Value = "0.0₃0678"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.331"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0411"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("E49").Value = "  +9.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.06%  "
